# Auto-generated PowerShell Word COM-interop script
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $ok = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for: $old"
    }
}

# Inserts Count new (initially blank, style-inheriting) paragraphs right after
# $d.Paragraphs($index), all anchored off the SAME original paragraph so none of
# them inherit formatting from one another. Returns the new paragraph count added.
function Insert-BlankParas($index, $count) {
    $anchor = $d.Paragraphs($index)
    for ($i = 0; $i -lt $count; $i++) {
        $anchor.Range.InsertParagraphAfter()
    }
}

# ---- Simple 1:1 text replacements (paragraph count unchanged) ----
Replace-Text "Distinguished Research and Data Analytics Leader with 21 years of experience directing groundbreaking applied research projects that have shaped policy, influenced elections, and transformed community development outcomes. Led multi-million dollar research initiatives serving thousands of analysts nationwide, with proven expertise in translating complex research insights for diverse stakeholders including elected officials, government agencies, and community organizations. Expert in research methodology design, statistical analysis, and community partnership development with extensive experience addressing systemic poverty challenges and delivering evidence-based solutions that drive meaningful social impact." "Distinguished Research & Data Analytics Professional with 21 years of expertise in survey methodology, consumer insights, voting behavior, and advanced data analysis. Proven track record in designing and implementing comprehensive research studies, managing cross-functional teams, and translating complex data into actionable intelligence. Expert in geospatial analysis, demographic segmentation, and consumer behavior modeling with experience serving major brands, organizations, and political candidates. Regular expert testimony and source on public opinion for journalists, with redistricting analysis used in court cases."
Replace-Text "Applied Research Leadership: Applied Research Project Management (Conception to Completion) • Research Methodology Design and Implementation • Cross-functional Team Leadership and Mentoring • Stakeholder Communication and Translation of Complex Findings • Evidence-Based Framework Development • Survey Methodology and Consumer Insights • Statistical Analysis and Data Validation" "Survey Methodology & Research Design: Survey Design and Questionnaire Development for Political and Market Research • Sampling Methodology and Statistical Analysis (R, SPSS, Stata, OSCAR) • Random Device Engagement (RDE), Text Message, Web Panel, and Live Telephone Calling • Focus Groups and Qualitative Research Methodologies • Meta-analytical Dataset Development for Longitudinal Analysis • Survey Instrument Standardization and Call Methods Optimization • Expert Testimony and Consultation on Research Methodology"
Replace-Text "Technical Proficiency: Programming: Python (Pandas, SciKit, TensorFlow, Django), R, SQL, Scala • Data Platforms: PostgreSQL, MySQL, Snowflake, Spark, MongoDB, Oracle • Analysis Tools: Excel (Advanced), Tableau, PowerBI, SPSS, SAS • Research Tools: Survey Design, Sampling Methodology, Statistical Modeling • Geospatial Analysis: ESRI ArcGIS, Quantum GIS, PostGIS, OSGeo" "Data Analysis & Visualization: Advanced Statistical Modeling and Analysis (Regression, Clustering, Segmentation) • Data Visualization: Tableau, PowerBI, Seaborn, Matplotlib, d3.js • Geospatial Analysis: ArcGIS, Quantum GIS, GRASS, OSGeo, PostGIS • Choropleths and Hexagonal Grid Maps for Demographic Visualization • Consumer Behavior Analysis and Market Segmentation • Machine Learning and Predictive Modeling for Targeting • Big Data Analytics: Spark/PySpark, Hadoop, Snowflake, dbt"
Replace-Text "Strategic Operations: Community Partnership Development • Government Relations and Policy Analysis • Multi-million Dollar Project Management • Performance Measurement and Evaluation • Data-Driven Decision Making for Social Impact • Public Systems Integration • Stakeholder Briefing and Expert Testimony" "Research Leadership & Client Management: Multi-million Dollar Research Project Management • Cross-functional Team Leadership (Teams of 7-11 professionals) • Client Relationship Management across Political, NGO, and Corporate Sectors • Stakeholder Briefing for Elected Officials and Senior Leadership • Court Case Analysis and Expert Testimony • Research Framework Development and Quality Control • Business Intelligence and Market Intelligence Delivery"
Replace-Text "• Conduct comprehensive quantitative and qualitative research studies using Python, R, SPSS, and Stata for political candidates and organizations" "• Conducted comprehensive quantitative and qualitative research studies for political candidates and major organizations, providing actionable consumer insights and market intelligence"
Replace-Text "• Architect cloud-based data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics" "• Designed and implemented advanced segmentation models using demographic, psychographic, and behavioral data to identify high-value targets"
Replace-Text "• Design scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets" "• Led multi-million dollar market research projects involving sensitive consumer data, ensuring compliance with privacy regulations while delivering actionable insights"
Replace-Text "• Develop custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering" "• Developed and deployed custom research software that processed billions of consumer records for pattern analysis, fraud detection and entity resolution"
Replace-Text "• Manage complex client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications" "• Built and maintained client relationships across diverse industries, consistently delivering insights that drove strategic decision-making"
Replace-Text "• Lead technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices" "• Architected and engineered cloud-based data warehouse solutions processing billions of records for electoral analytics and geospatial analysis"
Replace-Text "• Developed software solutions for political campaigns and advocacy groups" "• Developed software solutions for political campaigns and advocacy groups using modern web technologies"
Replace-Text "• Built web applications for voter engagement and campaign management" "• Built web applications for voter engagement and campaign management with real-time data integration"
Replace-Text "• Integrated third-party APIs and data sources for campaign tools" "• Integrated third-party APIs and data sources for campaign tools and voter database management"
Replace-Text "• Collaborated with political strategists to translate requirements into technical solutions" "• Collaborated with political strategists to translate requirements into technical solutions using agile methodologies"
Replace-Text "• Integrated technology solutions within organizational frameworks for social justice organizations" "• Integrated technology solutions within organizational frameworks for social justice organizations using open source technologies"
Replace-Text "• Developed data management systems for community organizing efforts" "• Developed data management systems for community organizing efforts with focus on accessibility and usability"
Replace-Text "• Provided technical training and support to nonprofit staff" "• Provided technical training and support to nonprofit staff on database management and data analysis tools"
Replace-Text "• Built custom applications for community engagement and advocacy" "• Built custom applications for community engagement and advocacy using web technologies and mobile platforms"

# ---- Sections needing in-place bullet replacement + new bullet insertion ----
# -- DATA PRODUCTS MANAGER bullet text replacements --
Replace-Text "• Conceived and developed framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES" "• Conceived and led implementation of comprehensive multi-tenant data warehouse integrating consumer demographic, economic, and behavioral data"
Replace-Text "• Built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions" "• Overhauled the organization's survey methodology and polling operations, significantly improving data accuracy and response rates"
Replace-Text "• Trained analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI" "• Managed and developed one of the in-house polling teams, focusing on Random Device Engagement (RDE), text message and web panel collected surveys, with live telephone calling and focus groups"
Replace-Text "• Wrote five-year strategic plans for developing data warehouse using Scala, PySpark, and Apache Spark that became basis of company's distinguishing products" "• Worked on standardizing questions, survey instruments and call methods, along with building a meta-analytical dataset for longitudinal analysis"
Replace-Text "• Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices" "• Managed a cross-functional team of eleven data engineers and analysts, establishing best practices for research methodology and data analysis"
# -- SENIOR ANALYST bullet text replacements --
Replace-Text "• Developed RACSO, a web application for pollsters to fully administer research including questionnaire creation, versioning, and reporting" "• Designed comprehensive survey instruments for specialized voting segments and niche markets"
Replace-Text "• Wrote RFP and analyzed bids from 1,200 vendors before selecting implementation partner" "• Developed sophisticated analytical products and reports that delivered actionable insights to clients"
Replace-Text "• Built prototype in R for comprehensive polling administration and sample file management" "• Co-developed RACSO web application to manage all aspects of survey operations, from instrument design to data collection and analysis"
Replace-Text "• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research" "• Wrote RFP and analyzed bids from 1,200 vendors before selecting implementation partner for RACSO platform"
# -- RESEARCH DIRECTOR bullet text replacements --
Replace-Text "• Managed critical research operations for political campaigns" "• Engineered FLEEM web application using Twilio's API to make thousands of simultaneous phone calls for IVR polls"
Replace-Text "• Conducted comprehensive polling and demographic analysis" "• Used FLEEM for early quantitative research in support of Senators Martin Heinrich and Elizabeth Warren"
Replace-Text "• Developed strategic recommendations based on data analysis" "• Led all aspects of survey design, implementation, data analysis, and reporting for major national studies"
Replace-Text "• Led research team in support of progressive political initiatives" "• Developed new statistical methods for boundary estimation techniques, enhancing geographic market segmentation capabilities"
# -- PROGRAMMER (Lake Research Partners) bullet text replacements --
Replace-Text "• Developed data analysis tools for political polling and research" "• Worked on all aspects of questionnaire design, sampling, reporting and analysis for political actors in Congressional, Senate and Presidential elections"
Replace-Text "• Built statistical models for voter behavior analysis" "• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party and affiliated actors"
Replace-Text "• Created data visualization tools for research presentations" "• Developed polling consortium database that later became the Polling Consortium Database at The Analyst Institute"
Replace-Text "• Supported senior researchers with technical analysis and reporting" "• Designed questionnaires and analyzed data for complex market research studies across diverse industries"
# -- FIELD DIRECTOR bullet text replacements --
Replace-Text "• Managed field operations for political campaigns and research projects" "• Administered all quantitative and qualitative research, ensuring that reporting was accurate and comprehensive"
Replace-Text "• Developed data collection and management systems for field work" "• Managed all aspects of survey fielding for a multi-million dollar research firm, including scheduling, oversight, sampling, and quality control"
Replace-Text "• Trained field staff on data collection protocols and quality control" "• Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings"
Replace-Text "• Analyzed field data to inform campaign strategy and research findings" "• Created custom reports and data visualizations based on specific client requirements"

# ---- Insert new bullets at fixed original-document positions, adjusting for
# paragraphs already inserted earlier in the document (running offset). ----
$offset = 0
# -- DATA PRODUCTS MANAGER: insert 1 new bullet(s) after original paragraph 24 --
$base = 24 + $offset
Insert-BlankParas $base 1
$d.Paragraphs($base + 1).Range.Text = "• Developed advanced data pipelines for machine learning applications that enhanced consumer segmentation and predictive modeling capabilities"
$offset = $offset + 1

# -- SENIOR ANALYST: insert 2 new bullet(s) after original paragraph 36 --
$base = 36 + $offset
Insert-BlankParas $base 2
$d.Paragraphs($base + 1).Range.Text = "• Introduced geospatial techniques to enhance market segmentation capabilities, providing clients with location-based consumer insights"
$d.Paragraphs($base + 2).Range.Text = "• Standardized reporting methodologies to improve clarity and impact of research findings"
$offset = $offset + 2

# -- RESEARCH DIRECTOR: insert 2 new bullet(s) after original paragraph 42 --
$base = 42 + $offset
Insert-BlankParas $base 2
$d.Paragraphs($base + 1).Range.Text = "• Created comprehensive data visualization solutions that improved clients' understanding of complex research findings"
$d.Paragraphs($base + 2).Range.Text = "• Provided tabular and graphical reporting with plans for interactive data exploration capabilities"
$offset = $offset + 2

# -- PROGRAMMER (Lake Research Partners): insert 2 new bullet(s) after original paragraph 60 --
$base = 60 + $offset
Insert-BlankParas $base 2
$d.Paragraphs($base + 1).Range.Text = "• Conducted statistical modeling and analysis to address multifaceted consumer behavior questions"
$d.Paragraphs($base + 2).Range.Text = "• Pioneered the integration of advanced mapping techniques into standard reports, including choropleths and hexagonal grid maps"
$offset = $offset + 2

# -- FIELD DIRECTOR: insert 2 new bullet(s) after original paragraph 66 --
$base = 66 + $offset
Insert-BlankParas $base 2
$d.Paragraphs($base + 1).Range.Text = "• Introduced mapping and geospatial analysis into standard reporting procedures, enhancing the value of research deliverables"
$d.Paragraphs($base + 2).Range.Text = "• Trained field staff on data collection protocols and quality control using standardized methodologies and best practices"
$offset = $offset + 2

# ---- KEY ACHIEVEMENTS AND IMPACT section ----
Replace-Text "Research Leadership and Policy Impact" "Survey Methodology & Research Innovation"
Replace-Text "• Regular expert testimony and consultation on research methodology for journalists, elected officials, and community leaders" "• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party, later becoming the Polling Consortium Database at The Analyst Institute"
Replace-Text "• Research analysis used in court cases addressing housing, redistricting, and community development with rigorous methodology" "• Developed RACSO platform for pollsters to fully administer research, analyzing bids from 1,200 vendors before selecting implementation partner"
Replace-Text "• Conceived and deployed cloud-based analytical software used by thousands of researchers nationwide for community-focused research" "• Engineered FLEEM system using Twilio API for thousands of simultaneous phone calls for IVR polls supporting Senators Martin Heinrich and Elizabeth Warren"
Replace-Text "• Developed research frameworks and methodologies that became industry standards for community development and policy analysis" "• Pioneered the integration of advanced mapping techniques into standard reports, including choropleths and hexagonal grid maps"

# Append two new Heading3 groups after original paragraph 72 (now shifted by $offset)
$base = 72 + $offset
Insert-BlankParas $base 5
$d.Paragraphs($base + 1).Range.Text = "Expert Testimony & Court Cases"
$d.Paragraphs($base + 1).Style = "Heading3"
$d.Paragraphs($base + 2).Range.Text = "• Regular expert testimony and source on public opinion for journalists, elected officials, and NGO leadership"
$d.Paragraphs($base + 3).Range.Text = "• Redistricting analysis used in court cases with rigorous methodology and expert testimony"
$d.Paragraphs($base + 4).Range.Text = "• Research analysis used in court cases addressing housing, redistricting, and community development"
$d.Paragraphs($base + 5).Range.Text = "• Provided expert consultation on research methodology for diverse stakeholder groups"
$base = $base + 5

Insert-BlankParas $base 5
$d.Paragraphs($base + 1).Range.Text = "Data Infrastructure & Analytics"
$d.Paragraphs($base + 1).Style = "Heading3"
$d.Paragraphs($base + 2).Range.Text = "• Conceived, architected, engineered and deployed cloud-based redistricting software used by thousands of analysts nationwide"
$d.Paragraphs($base + 3).Range.Text = "• Designed, architected and created multi-tenant data warehouse tracking decades of political, geographical, econometric change"
$d.Paragraphs($base + 4).Range.Text = "• Led multi-million dollar market research projects involving sensitive consumer data with privacy compliance"
$d.Paragraphs($base + 5).Range.Text = "• Developed research frameworks and methodologies that became industry standards for community development and policy analysis"
$base = $base + 5

Write-Host "DONE"